# Insert 16 new rows before row 102 (pushing existing rows 102-103 down to 118-119)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("102:117").Insert()

$ws.Range("A102").Value = "HSBCDOL"
$ws.Range("B102").Value = "BNF"
$ws.Range("C102").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D102").Value = "Mercado de Deuda Extranjero"
$ws.Range("E102").Value = "Renta Fija Internacional"
$ws.Range("A103").Value = "NTEDLS+"
$ws.Range("B103").Value = "FF"
$ws.Range("C103").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D103").Value = "Mercado de Deuda Extranjero"
$ws.Range("E103").Value = "Renta Fija Internacional"
$ws.Range("A104").Value = "PRGLOB"
$ws.Range("B104").Value = "FFR"
$ws.Range("C104").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D104").Value = "Mercado de Deuda Extranjero"
$ws.Range("E104").Value = "Renta Fija Internacional"
$ws.Range("A105").Value = "SCOTDL+"
$ws.Range("B105").Value = "C1E"
$ws.Range("C105").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D105").Value = "Mercado de Deuda Extranjero"
$ws.Range("E105").Value = "Renta Fija Internacional"
$ws.Range("A106").Value = "SCOTDOL"
$ws.Range("B106").Value = "C1E"
$ws.Range("C106").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D106").Value = "Mercado de Deuda Extranjero"
$ws.Range("E106").Value = "Renta Fija Internacional"
$ws.Range("A107").Value = "SURUSD"
$ws.Range("B107").Value = "BOE"
$ws.Range("C107").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D107").Value = "Mercado de Deuda Extranjero"
$ws.Range("E107").Value = "Renta Fija Internacional"
$ws.Range("A108").Value = "PRINHYD"
$ws.Range("B108").Value = "FFX"
$ws.Range("C108").Value = "Acciones de Sociedades de Inversion de Instrumentos de Deuda"
$ws.Range("D108").Value = "Mercado de Deuda Extranjero"
$ws.Range("E108").Value = "Renta Fija Internacional"
$ws.Range("A109").Value = "NTEINT+"
$ws.Range("B109").Value = "FF"
$ws.Range("C109").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D109").Value = "Mercado de Capitales Extranjero"
$ws.Range("E109").Value = "Renta Variable Internacional"
$ws.Range("A110").Value = "NTEUSA"
$ws.Range("B110").Value = "FF"
$ws.Range("C110").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D110").Value = "Mercado de Capitales Extranjero"
$ws.Range("E110").Value = "Renta Variable Internacional"
$ws.Range("A111").Value = "NTEUSA+"
$ws.Range("B111").Value = "FF"
$ws.Range("C111").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D111").Value = "Mercado de Capitales Extranjero"
$ws.Range("E111").Value = "Renta Variable Internacional"
$ws.Range("A112").Value = "SCOTGLO"
$ws.Range("B112").Value = "C1E"
$ws.Range("C112").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D112").Value = "Mercado de Capitales Extranjero"
$ws.Range("E112").Value = "Renta Variable Internacional"
$ws.Range("A113").Value = "SURGLOB"
$ws.Range("B113").Value = "BOE0"
$ws.Range("C113").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D113").Value = "Mercado de Capitales Extranjero"
$ws.Range("E113").Value = "Renta Variable Internacional"
$ws.Range("A114").Value = "BLKINT1"
$ws.Range("B114").Value = "M0-A"
$ws.Range("C114").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D114").Value = "Mercado de Capitales Extranjero"
$ws.Range("E114").Value = "Renta Variable Internacional"
$ws.Range("A115").Value = "FT-GLOB"
$ws.Range("B115").Value = "BE"
$ws.Range("C115").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D115").Value = "Mercado de Capitales Extranjero"
$ws.Range("E115").Value = "Renta Variable Internacional"
$ws.Range("A116").Value = "PEMERGE"
$ws.Range("B116").Value = "FFR"
$ws.Range("C116").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D116").Value = "Mercado de Capitales Extranjero"
$ws.Range("E116").Value = "Renta Variable Internacional"
$ws.Range("A117").Value = "PRINFUS"
$ws.Range("B117").Value = "FFX"
$ws.Range("C117").Value = "Acciones de Sociedades de Inversion"
$ws.Range("D117").Value = "Mercado de Capitales Extranjero"
$ws.Range("E117").Value = "Renta Variable Internacional"

# Update the selection/view to match the new scroll position
$ws.Range("A99").Select()
